$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DocumentManifest")

# Column F ("Comments") gains FHIR Reference-type annotations for a handful
# of rows that previously had an empty Comments cell.
$ws.Range("F7").Value  = "Reference(Patient | Practitioner | Group | Device)"
$ws.Range("F9").Value  = "Reference(Practitioner | Organization | Device | Patient | RelatedPerson)"
$ws.Range("F10").Value = "Reference(Patient | Practitioner | RelatedPerson | Organization)"
$ws.Range("F16").Value = "Reference(Any)"
$ws.Range("F19").Value = "Reference(Any)"
